# Apply the LinuxForHealth re-brand / version bump edit to the
# StructureDefinition-employee-hire-date workbook.

$wb = $excel.ActiveWorkbook

# ---- "Metadata" worksheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# URL (row 2) - ibm.com -> linuxforhealth.org
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-hire-date"

# Version (row 3) - 7.0.0 -> 8.0.0
$meta.Range("B3").Value = "8.0.0"

# Date (row 8) - publication timestamp refresh
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher (row 9) - Alvearie Team -> LinuxForHealth Team
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---- "Elements" worksheet --------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# The combined ele-1/ext-1 invariant text used to be (incorrectly) duplicated
# on the root "Extension" row (row 2, column AI "Constraint(s)"). It belongs
# solely on the "Extension.extension" row (row 4), which already carries it,
# so clear the stray copy on row 2.
$elements.Range("AI2").Value = ""

# "Extension.url" row (row 5) fixes its value to the extension's own
# canonical URL, so it must track the same ibm.com -> linuxforhealth.org
# rename applied above.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-hire-date"
